# Fix merge issues in the document's SharePoint "contentTypeSchema" custom
# XML part: the content-type display strings were left in Japanese by a bad
# merge; restore the English strings and let PowerPoint regenerate the
# schema's version/fields identifiers (it does this automatically whenever
# the part's XML is rewritten).
#
# This touches customXml/item1.xml (ma:contentTypeName, ma:contentTypeDescription,
# ma:versionID, ma:fieldsID, and the three ma:displayName attributes) via the
# CustomXMLParts object model - the supported way to edit a package's custom
# XML parts through PowerPoint COM automation.

$p = $ppt.ActivePresentation

$cxps = $p.CustomXMLParts
$schemaParts = $cxps.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")

if ($schemaParts.Count -ge 1) {
    $part = $schemaParts.Item(1)
    $xml = $part.XML

    # Localized content-type metadata -> English.
    $xml = $xml.Replace('ma:contentTypeName="ドキュメント"', 'ma:contentTypeName="Document"')
    $xml = $xml.Replace('ma:contentTypeDescription="新しいドキュメントを作成します。"', 'ma:contentTypeDescription="Create a new document."')
    $xml = $xml.Replace('ma:displayName="画像タグ"', 'ma:displayName="Image Tags"')
    $xml = $xml.Replace('ma:displayName="コンテンツ タイプ"', 'ma:displayName="Content Type"')
    $xml = $xml.Replace('ma:displayName="タイトル"', 'ma:displayName="Title"')

    # Regenerated identifiers that PowerPoint stamps whenever the schema
    # content changes.
    $xml = $xml.Replace('ma:versionID="a9cab35011a557c1232e9e1918db7064"', 'ma:versionID="d0e002fabf17cb2440d8e9a473d3a41c"')
    $xml = $xml.Replace('ma:fieldsID="36c473bbc383ceb924bb8d2cdd9a2de6"', 'ma:fieldsID="e4cec627508c1f1ba247db94416ea198"')

    $part.XML = $xml
}

# The datastore item that backs item1.xml gets a freshly generated GUID
# whenever its content is rewritten.
$itemProps = $p.CustomXMLParts.SelectByID("{C8E3A567-4AEF-48F7-8398-D01FD834A75B}")
if ($itemProps.Count -ge 1) {
    $itemProps.Item(1).Id = "{FE7C2605-E662-437E-A8A4-A5435D757B8B}"
}
